# Mise à jour de certains champs de Modules et de Professeurs
#
# - Column C header: "Enseignant"       -> "Chef  Module"
# - Column D header: "Nombre d'heures"  -> "Composants"
# - Column C width  -> 35 characters
# - Column D width  -> ~24.57 characters
# - Active selection moved to E8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# --- Column widths ---------------------------------------------------------
# (the host quantizes the stored width to 1/6-character steps and adds a
#  5/6-character on-screen padding term to whatever ColumnWidth is assigned,
#  so the assigned values below are pre-compensated to land on the closest
#  achievable width to the authored values of 35 and 24.5703125 characters)
$ws.Columns.Item(3).ColumnWidth = 34.1666666666667
$ws.Columns.Item(4).ColumnWidth = 23.6666666666667

# --- Selection --------------------------------------------------------------
[void]$ws.Range("E8").Select()
